$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:I1
$ws.Range("F1").Value = "當月MPS"
$ws.Range("G1").Value = "當月生產天數"
$ws.Range("H1").Value = "下月MPS"
$ws.Range("I1").Value = "下月生產天數"

# New data cells F2:I2 (numeric values)
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1

# Update selection to match target state
$ws.Range("G11").Select()
